$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transit Mode")

# Insert a new row before row 67 (which currently holds mode 110 / MUNI Metro),
# shifting subsequent rows down, so the new row becomes row 67.
$ws.Rows.Item(67).Insert()

# Fill in the new row 67: Transit mode 107 "Treasure Island Ferry", Ferry category, WETA agency, FB NTD mode
$ws.Cells.Item(67, 1).Value = 107
$ws.Cells.Item(67, 2).Value = "Treasure Island Ferry"
$ws.Cells.Item(67, 3).Value = "Ferry"
$ws.Cells.Item(67, 4).Value = "WETA"
$ws.Cells.Item(67, 5).Value = "FB"

# Copy formatting from the row above (row 66), restricted to columns A:E, into new row 67
$ws.Range("A66:E66").Copy()
$ws.Range("A67:E67").PasteSpecial(-4122)  # xlPasteFormats

# Restore the values (paste of formats shouldn't touch them, but ensure they are correct)
$ws.Cells.Item(67, 1).Value = 107
$ws.Cells.Item(67, 2).Value = "Treasure Island Ferry"
$ws.Cells.Item(67, 3).Value = "Ferry"
$ws.Cells.Item(67, 4).Value = "WETA"
$ws.Cells.Item(67, 5).Value = "FB"

# Update row 65 (South San Francisco Ferry, mode 105): fill in agency (WETA) and NTD mode (FB)
$ws.Cells.Item(65, 4).Value = "WETA"
$ws.Cells.Item(65, 5).Value = "FB"

$excel.CutCopyMode = 0
